$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 43.21270066666667
$ws.Range("H2").Value = 129.638102
$ws.Range("I2").Value = 0.1487696778665633
$ws.Range("J2").Value = 0.1487696778665633
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 56.98117766666667
$ws.Range("N2").Value = 170.943533
$ws.Range("O2").Value = 0.952030123851636
$ws.Range("P2").Value = 0.9520301238516359
$ws.Range("Q2").Value = 2462.310574143818
$ws.Range("R2").Value = 22160.79516729437
$ws.Range("S2").Value = 0.1416332148446723
$ws.Range("T2").Value = 0.1416332148446722

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 43.21270066666667
$ws.Range("H3").Value = 129.638102
$ws.Range("I3").Value = 0.1487696778665633
$ws.Range("J3").Value = 0.1487696778665633
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.516719
$ws.Range("N3").Value = 7.550157
$ws.Range("O3").Value = 0.04204883786863874
$ws.Range("P3").Value = 0.04204883786863874
$ws.Range("Q3").Value = 108.7542248091127
$ws.Range("R3").Value = 978.7880232820141
$ws.Range("S3").Value = 0.006255592064380734
$ws.Range("T3").Value = 0.006255592064380733

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 43.21270066666667
$ws.Range("H4").Value = 129.638102
$ws.Range("I4").Value = 0.1487696778665633
$ws.Range("J4").Value = 0.1487696778665633
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3543876666666666
$ws.Range("N4").Value = 1.063163
$ws.Range("O4").Value = 0.005921038279725251
$ws.Range("P4").Value = 0.005921038279725251
$ws.Range("Q4").Value = 15.31404815962511
$ws.Range("R4").Value = 137.826433436626
$ws.Range("S4").Value = 0.0008808709575103158
$ws.Range("T4").Value = 0.0008808709575103156

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 237.7114053333333
$ws.Range("H5").Value = 713.134216
$ws.Range("I5").Value = 0.8183762794517323
$ws.Range("J5").Value = 0.8183762794517323
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 56.98117766666667
$ws.Range("N5").Value = 170.943533
$ws.Range("O5").Value = 0.952030123851636
$ws.Range("P5").Value = 0.9520301238516359
$ws.Range("Q5").Value = 13545.07582069168
$ws.Range("R5").Value = 121905.6823862251
$ws.Range("S5").Value = 0.7791188706836738
$ws.Range("T5").Value = 0.7791188706836737

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 237.7114053333333
$ws.Range("H6").Value = 713.134216
$ws.Range("I6").Value = 0.8183762794517323
$ws.Range("J6").Value = 0.8183762794517323
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.516719
$ws.Range("N6").Value = 7.550157
$ws.Range("O6").Value = 0.04204883786863874
$ws.Range("P6").Value = 0.04204883786863874
$ws.Range("Q6").Value = 598.2528103191014
$ws.Range("R6").Value = 5384.275292871913
$ws.Range("S6").Value = 0.03441177149020568
$ws.Range("T6").Value = 0.03441177149020568

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 237.7114053333333
$ws.Range("H7").Value = 713.134216
$ws.Range("I7").Value = 0.8183762794517323
$ws.Range("J7").Value = 0.8183762794517323
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.3543876666666666
$ws.Range("N7").Value = 1.063163
$ws.Range("O7").Value = 0.005921038279725251
$ws.Range("P7").Value = 0.005921038279725251
$ws.Range("Q7").Value = 84.24199027613422
$ws.Range("R7").Value = 758.1779124852079
$ws.Range("S7").Value = 0.004845637277852836
$ws.Range("T7").Value = 0.004845637277852836

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.543019333333334
$ws.Range("H8").Value = 28.629058
$ws.Range("I8").Value = 0.03285404268170446
$ws.Range("J8").Value = 0.03285404268170446
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 56.98117766666667
$ws.Range("N8").Value = 170.943533
$ws.Range("O8").Value = 0.952030123851636
$ws.Range("P8").Value = 0.9520301238516359
$ws.Range("Q8").Value = 543.7724801091016
$ws.Range("R8").Value = 4893.952320981914
$ws.Range("S8").Value = 0.03127803832329003
$ws.Range("T8").Value = 0.03127803832329003

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.543019333333334
$ws.Range("H9").Value = 28.629058
$ws.Range("I9").Value = 0.03285404268170446
$ws.Range("J9").Value = 0.03285404268170446
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.516719
$ws.Range("N9").Value = 7.550157
$ws.Range("O9").Value = 0.04204883786863874
$ws.Range("P9").Value = 0.04204883786863874
$ws.Range("Q9").Value = 24.01709807356734
$ws.Range("R9").Value = 216.153882662106
$ws.Range("S9").Value = 0.001381474314052328
$ws.Range("T9").Value = 0.001381474314052328

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.543019333333334
$ws.Range("H10").Value = 28.629058
$ws.Range("I10").Value = 0.03285404268170446
$ws.Range("J10").Value = 0.03285404268170446
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3543876666666666
$ws.Range("N10").Value = 1.063163
$ws.Range("O10").Value = 0.005921038279725251
$ws.Range("P10").Value = 0.005921038279725251
$ws.Range("Q10").Value = 3.381928354494888
$ws.Range("R10").Value = 30.437355190454
$ws.Range("S10").Value = 0.0001945300443620994
$ws.Range("T10").Value = 0.0001945300443620994

